# Update Pooh Points site
# Applies the 2026-01-24 refresh to the "Players" sheet:
#  - bumps the live game clock from "3:04 - 1st Half" to "2:44 - 1st Half"
#    for every player still mid-game (status column G)
#  - bumps O47 (min) from 10 to 11
#  - re-shuffles rows 99-104 (Undrafted / LSU@ARK & UGA@TEX & VAN@MSST & SC@TA&M
#    games) to reflect each player's latest box score

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# --- Simple status-clock refresh (G column), rows where the game is still live ---
$statusRows = @(8, 10, 15, 18, 25, 32, 39, 43, 47, 66, 67, 91, 92, 109)
foreach ($r in $statusRows) {
    $ws.Cells.Item($r, 7).Value = "2:44 - 1st Half"
}

# O47 (min) ticks up from 10 to 11
$ws.Cells.Item(47, 15).Value = 11

# --- Rows 99-104: refreshed box scores for the Undrafted LSU@ARK / UGA@TEX / VAN@MSST / SC@TA&M players ---

# Row 99: Rashad King (LSU)
$ws.Cells.Item(99, 4).Value = "Rashad King"
$ws.Cells.Item(99, 5).Value = "LSU"
$ws.Cells.Item(99, 6).Value = "LSU@ARK"
$ws.Cells.Item(99, 7).Value = "2:44 - 1st Half"
$ws.Cells.Item(99, 8).Value = 4
$ws.Cells.Item(99, 9).Value = 4
$ws.Cells.Item(99, 10).Value = 1
$ws.Cells.Item(99, 11).Value = 1
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0
$ws.Cells.Item(99, 14).Value = 0
$ws.Cells.Item(99, 15).Value = 9

# Row 100: Simeon Wilcher (TEX)
$ws.Cells.Item(100, 4).Value = "Simeon Wilcher"
$ws.Cells.Item(100, 5).Value = "TEX"
$ws.Cells.Item(100, 6).Value = "UGA@TEX"
$ws.Cells.Item(100, 7).Value = "Final"
$ws.Cells.Item(100, 8).Value = 4
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 4
$ws.Cells.Item(100, 11).Value = 2
$ws.Cells.Item(100, 12).Value = 2
$ws.Cells.Item(100, 13).Value = 0
$ws.Cells.Item(100, 14).Value = 0
$ws.Cells.Item(100, 15).Value = 15

# Row 101: Dellquan Warren (MSST)
$ws.Cells.Item(101, 4).Value = "Dellquan Warren"
$ws.Cells.Item(101, 5).Value = "MSST"
$ws.Cells.Item(101, 6).Value = "VAN@MSST"
$ws.Cells.Item(101, 7).Value = "Final"
$ws.Cells.Item(101, 8).Value = 3
$ws.Cells.Item(101, 9).Value = 2
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 2
$ws.Cells.Item(101, 12).Value = 1
$ws.Cells.Item(101, 13).Value = 0
$ws.Cells.Item(101, 14).Value = 2
$ws.Cells.Item(101, 15).Value = 9

# Row 102: Josh Holloway (TA&M)
$ws.Cells.Item(102, 4).Value = "Josh Holloway"
$ws.Cells.Item(102, 5).Value = "TA&M"
$ws.Cells.Item(102, 6).Value = "SC@TA&M"
$ws.Cells.Item(102, 7).Value = "Final"
$ws.Cells.Item(102, 8).Value = 3
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 4
$ws.Cells.Item(102, 12).Value = 1
$ws.Cells.Item(102, 13).Value = 0
$ws.Cells.Item(102, 14).Value = 1
$ws.Cells.Item(102, 15).Value = 12

# Row 103: Kareem Stagg (UGA)
$ws.Cells.Item(103, 4).Value = "Kareem Stagg"
$ws.Cells.Item(103, 5).Value = "UGA"
$ws.Cells.Item(103, 6).Value = "UGA@TEX"
$ws.Cells.Item(103, 7).Value = "Final"
$ws.Cells.Item(103, 8).Value = 3
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = 4
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = 0
$ws.Cells.Item(103, 14).Value = 0
$ws.Cells.Item(103, 15).Value = 13

# Row 104: Malique Ewin (ARK)
$ws.Cells.Item(104, 4).Value = "Malique Ewin"
$ws.Cells.Item(104, 5).Value = "ARK"
$ws.Cells.Item(104, 6).Value = "LSU@ARK"
$ws.Cells.Item(104, 7).Value = "2:44 - 1st Half"
$ws.Cells.Item(104, 8).Value = 3
$ws.Cells.Item(104, 9).Value = 6
$ws.Cells.Item(104, 10).Value = 1
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 1
$ws.Cells.Item(104, 13).Value = 1
$ws.Cells.Item(104, 14).Value = 1
$ws.Cells.Item(104, 15).Value = 10
